# ---------------------------------------------------------------------------
# silver_fore.xlsx: insert two new leading "Unnamed: 0.x" index columns
# (shifting the existing B:F block to D:H) and append 10 more data rows
# (22-31) that continue the same index / SILVER_FOR staircase pattern.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell that already carries the workbook's header/index style (bold font,
# thin box border, centered-top alignment) so we can stamp it onto the new
# cells instead of re-describing the format by hand.
$styleSrc = $ws.Cells.Item(2, 1)

# --- 1. Header row -----------------------------------------------------
# Old layout:  B=Unnamed:0.2 C=Unnamed:0.1 D=Unnamed:0 E="   SILVER_FOR" F=SILVER_FOR
# New layout:  B=Unnamed:0.4 C=Unnamed:0.3 D=Unnamed:0.2 E=Unnamed:0.1 F=Unnamed:0 G="   SILVER_FOR" H=SILVER_FOR
$ws.Cells.Item(1, 2).Value2 = "Unnamed: 0.4"
$ws.Cells.Item(1, 3).Value2 = "Unnamed: 0.3"
$ws.Cells.Item(1, 4).Value2 = "Unnamed: 0.2"
$ws.Cells.Item(1, 5).Value2 = "Unnamed: 0.1"
$ws.Cells.Item(1, 6).Value2 = "Unnamed: 0"
$ws.Cells.Item(1, 7).Value2 = "   SILVER_FOR"
$ws.Cells.Item(1, 8).Value2 = "SILVER_FOR"

# B1:F1 already carried the header style before this edit, so they keep it
# automatically; G1/H1 are brand-new cells and need it applied explicitly.
$styleSrc.Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(1, 8).PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Column A (row index, style 1) extends from 21 rows to 31 rows --
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}
$styleSrc.Copy()
$ws.Range($ws.Cells.Item(22, 1), $ws.Cells.Item(31, 1)).PasteSpecial(-4122)

# --- 3. Columns B-F: the same 0..N "Unnamed" staircase, one column     --
#        shorter than the previous, each new column inheriting the      --
#        longest run (B keeps the full 0..24 index that used to live in --
#        old column B, and so on down to F = old column D's 0..4 run).  --
for ($r = 2; $r -le 26; $r++) { $ws.Cells.Item($r, 2).Value2 = $r - 2 }   # B2:B26  -> 0..24
for ($r = 2; $r -le 21; $r++) { $ws.Cells.Item($r, 3).Value2 = $r - 2 }   # C2:C21  -> 0..19
for ($r = 2; $r -le 16; $r++) { $ws.Cells.Item($r, 4).Value2 = $r - 2 }   # D2:D16  -> 0..14
for ($r = 2; $r -le 11; $r++) { $ws.Cells.Item($r, 5).Value2 = $r - 2 }   # E2:E11  -> 0..9
for ($r = 2; $r -le 6;  $r++) { $ws.Cells.Item($r, 6).Value2 = $r - 2 }   # F2:F6   -> 0..4

# --- 4. Column G: the original "   SILVER_FOR" values (rows 2-6 only) --
$gValues = @(27.44466, 26.933548, 27.173124, 26.856646, 26.422922)
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value2 = $gValues[$i]
}

# --- 5. Column H: the SILVER_FOR value series, rows 7-31 ---------------
$hValues = @(
    30.93059290717292, 30.6353459147918, 30.30023174736436, 29.53532024840348, 29.15916464141611,
    31.76244298992617, 32.23934168175691, 32.5524573182023, 32.38868114596181, 32.10415671664822,
    28.81685345771996, 28.51061683625062, 27.67123246313037, 27.34084839285276, 26.94339908891732,
    29.09709106441289, 28.96379896414152, 28.94952882821161, 29.52255201468267, 28.89384280964123,
    28.95396845664823, 28.94541683747048, 28.70188153581586, 27.87425082086361, 28.09067066164505
)
for ($i = 0; $i -lt $hValues.Length; $i++) {
    $ws.Cells.Item(7 + $i, 8).Value2 = $hValues[$i]
}

Write-Host "done"
